$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns (rows 2-51) remain text, not auto-converted to numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.973.11"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "1.650.98"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "309.73"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "0.3897"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "52.40"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "1.349"
$ws.Range("E10").Value = "  -4.20%  "
$ws.Range("D11").Value = "0.9998"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "0.08448"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "23.81"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "7.072"
$ws.Range("E14").Value = "  -3.60%  "
$ws.Range("D15").Value = "7.997"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "0.00001310"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").Value = "1.642.87"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "0.06997"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "19.68"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").Value = "6.978"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "13.80"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "24.001.09"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "2.444"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "2.972"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "22.09"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "152.38"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").Value = "5.409"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "138.09"
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("D31").Value = "7.935"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "2.511"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "1.828.98"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("D34").Value = "1.025"
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("D35").Value = "0.08092"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").Value = "6.731"
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("D37").Value = "0.02931"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2679"
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "10.72"
$ws.Range("E39").Value = "  -5.05%  "
$ws.Range("D40").Value = "0.09112"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("D41").Value = "0.7598"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").Value = "13.40"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "1.421"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").Value = "16.29"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").Value = "0.6961"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").Value = "2.466"
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("D47").Value = "4.101"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "0.9998"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "0.08333"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").Value = "134.79"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "1.223"
$ws.Range("E51").Value = "  -3.89%  "

# Restore default style (remove the temporary Text number format) so only
# cell values/content differ from the original, matching the authored diff
$ws.Range("D2:E51").Style = "Normal"
